$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 2461
$ws.Range("F5").Value = 1652
$ws.Range("F6").Value = 100
$ws.Range("F7").Value = 312
$ws.Range("F8").Value = 612
$ws.Range("F9").Value = 3497
$ws.Range("F10").Value = 918
$ws.Range("F11").Value = 1147
$ws.Range("F12").Value = 1565
$ws.Range("F13").Value = 28
$ws.Range("F14").Value = 881
$ws.Range("F16").Value = 1239
$ws.Range("F17").Value = 1780
$ws.Range("F19").Value = 443
$ws.Range("F20").Value = 1537
$ws.Range("F21").Value = 1070
$ws.Range("F22").Value = 2082
$ws.Range("F24").Value = 4230
$ws.Range("F25").Value = 44
$ws.Range("F26").Value = 2692
$ws.Range("F27").Value = 1207

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F23").Value = 108
$ws.Range("F24").Value = 5
$ws.Range("F28").Value = 11
$ws.Range("F36").Value = 423
$ws.Range("F41").Value = 14

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2537
$ws.Range("F5").Value = 2547
$ws.Range("F10").Value = 384
$ws.Range("F11").Value = 2952
$ws.Range("F12").Value = 458
$ws.Range("F13").Value = 790
$ws.Range("F14").Value = 190

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2537
$ws.Range("F5").Value = 2461
$ws.Range("F8").Value = 2952
$ws.Range("F9").Value = 790
$ws.Range("F12").Value = 100
$ws.Range("F13").Value = 312
$ws.Range("F14").Value = 612
$ws.Range("F16").Value = 918
$ws.Range("F17").Value = 1147
$ws.Range("F19").Value = 28
$ws.Range("F20").Value = 881
$ws.Range("F23").Value = 1239
$ws.Range("F30").Value = 1780
$ws.Range("F31").Value = 443
$ws.Range("F33").Value = 1537
$ws.Range("F35").Value = 108
$ws.Range("F36").Value = 108
$ws.Range("F38").Value = 1070
$ws.Range("F40").Value = 2082
$ws.Range("F41").Value = 11
$ws.Range("F44").Value = 4231
$ws.Range("F45").Value = 423
$ws.Range("F46").Value = 2692
